$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47; this shifts existing rows 47..125 down to 48..126
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new data point.
$ws.Cells.Item(47, 1).Value = 7
$ws.Cells.Item(47, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(47, 3).Value = "Ñuble"
$ws.Cells.Item(47, 4).Value = 44477
$ws.Cells.Item(47, 5).Value = 16
$ws.Cells.Item(47, 6).Value = 100112017
$ws.Cells.Item(47, 7).Value = "Apio"
$ws.Cells.Item(47, 8).Value = "Americana (o)"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 160
$ws.Cells.Item(47, 11).Value = 8000
$ws.Cells.Item(47, 12).Value = 9000
$ws.Cells.Item(47, 13).Value = 8500
$ws.Cells.Item(47, 14).Value = "$/docena de matas"
$ws.Cells.Item(47, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(47, 16).Value = 1417
$ws.Cells.Item(47, 17).Value = 6
$ws.Cells.Item(47, 18).Value = "Hortaliza"
